$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column BO (67), rows 3 through 29, with a single space value.
# This mirrors the author's action of adding a "refresh" column of blank/space
# marker cells below the header/team rows, extending the used range from
# A1:CA14 down to A1:CA29.
for ($r = 3; $r -le 29; $r++) {
    $ws.Cells.Item($r, 67).Value = " "
}

# Select the last cell of the new range, matching the final selection state.
$ws.Range("BO29").Select()
